$wb = $excel.ActiveWorkbook

# ---------- Sheet "hpi" ----------
$ws = $wb.Worksheets.Item("hpi")
$ws.Range("B2").Value = "Difficulty swallowing solids is a common symptom of food impaction, indicating a physical obstruction in the esophagus."
$ws.Range("D2").Value = "If antacids are no longer providing relief, it may suggest that the symptoms are not due to food impaction but rather to another underlying condition such as esophagitis or GERD."
$ws.Range("B3").Value = "Odynophagia is often associated with food impaction, as the presence of an obstructing food bolus can cause pain during swallowing."
$ws.Range("B4").Value = "Reflux can be a contributing factor to food impaction, as it may indicate esophageal motility issues or structural abnormalities."
$ws.Range("D4").Value = "Chronic reflux symptoms suggest a different underlying pathology, such as GERD, rather than food impaction, which is usually characterized by more acute episodes."
$ws.Range("B6").Value = "A prolonged duration of symptoms related to food getting stuck is highly indicative of food impaction, as it suggests a chronic issue."
$ws.Range("C6").Value = "Current heartburn is present."
$ws.Range("D6").Value = "Current heartburn is more indicative of acid reflux or esophagitis rather than food impaction, which typically presents with more specific symptoms related to swallowing."

# ---------- Sheet "hist" ----------
$ws = $wb.Worksheets.Item("hist")
$ws.Range("D2").Value = "Alcohol use disorder is a risk factor for esophageal issues; its absence suggests a lower likelihood of food impaction."
$ws.Range("D3").Value = "Nicotine dependence can lead to esophageal motility disorders; its absence may indicate a lower risk of food impaction."
$ws.Range("D4").Value = "Radiation treatment can lead to esophageal strictures, which increase the risk of food impaction; its absence suggests a lower risk."
$ws.Range("D5").Value = "Obesity is a risk factor for esophageal problems; its absence may indicate a lower likelihood of food impaction."
$ws.Range("D6").Value = "Hypertension can be associated with other comorbidities affecting esophageal motility; its absence may suggest a lower risk of food impaction."

# ---------- Sheet "soc" ----------
$ws = $wb.Worksheets.Item("soc")
$ws.Range("B2").Value = "Social stress can lead to behaviors that increase the risk of food impaction, such as neglecting proper eating habits or consuming inappropriate food types."
$ws.Range("D2").Value = "Absence of alcohol use is a strong indicator against food impaction, as alcohol can contribute to gastrointestinal issues that may lead to such conditions."
$ws.Range("B3").Value = "Tobacco use can contribute to gastrointestinal issues, which may increase the likelihood of food impaction."
$ws.Range("D3").Value = "Current tobacco use is often associated with gastrointestinal problems; its absence suggests a lower likelihood of food impaction."
$ws.Range("B5").Value = "While this finding is not directly supportive, the absence of autoimmune conditions may suggest a lower risk of related complications that could lead to food impaction."
$ws.Range("B6").Value = "Similar to the previous finding, the absence of cancer may indicate a lower risk of malignancies that could obstruct the gastrointestinal tract, thus indirectly supporting the diagnosis of food impaction."

# ---------- Sheet "obj" ----------
$ws = $wb.Worksheets.Item("obj")
$ws.Range("D2").Value = "A hoarse voice can indicate esophageal obstruction or irritation, which is often present in cases of food impaction."
$ws.Range("D3").Value = "Absence of cough may suggest that there is no acute respiratory issue, which can sometimes accompany food impaction."
$ws.Range("B4").Value = "Absence of cough may suggest that there is no acute respiratory issue, which can sometimes accompany food impaction."
$ws.Range("C4").Value = "Neck masses or fullness observed is absent."
$ws.Range("D4").Value = "Absence of neck masses or fullness may indicate that there is no obstruction or swelling in the area that could lead to food impaction."
$ws.Range("A5").Value = "Epigastric pain on palpation is absent."
$ws.Range("B5").Value = "Absence of epigastric pain may indicate that there is no acute gastrointestinal distress, which can sometimes be confused with food impaction."
$ws.Range("C5").Value = "Joint swelling observed is absent."
$ws.Range("D5").Value = "Absence of joint swelling may suggest that there is no systemic inflammatory condition that could contribute to esophageal issues."
$ws.Range("A6").Value = "Halitosis observed is absent."
$ws.Range("B6").Value = "Absence of halitosis may suggest that there is no prolonged food retention in the esophagus, which is often associated with food impaction."
$ws.Range("C6").Value = "Rheumatoid nodules is absent."
$ws.Range("D6").Value = "Absence of rheumatoid nodules may indicate that there is no underlying autoimmune condition that could affect the esophagus and lead to food impaction."
